# Apply updated cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '34.341.35'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.804.81'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +1.03%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '227.44'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.56%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.574'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.66%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '36.20'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +11.11%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.302'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +2.62%  '
$ws.Range('E10').Value = '  +0.89%  '
$ws.Range('E11').Value = '  +1.99%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '2.065.49'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '11.80'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +7.37%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.804.33'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.08%  '
$ws.Range('E15').Value = '  +2.15%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.52'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +5.77%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '34.343.59'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '69.14'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.38%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '245.73'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.56%  '
$ws.Range('E20').Value = '  +0.18%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.60'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +3.81%  '
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('E23').Value = '  +1.05%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '172.04'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +3.73%  '
$ws.Range('E25').Value = '  +2.64%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.99'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +9.73%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.90'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.62%  '
$ws.Range('E28').Value = '  +2.90%  '
$ws.Range('E29').Value = '  -0.25%  '
$ws.Range('E30').Value = '  +1.42%  '
$ws.Range('E31').Value = '  +1.52%  '
$ws.Range('E32').Value = '  +1.17%  '
$ws.Range('E33').Value = '  +0.87%  '
$ws.Range('E34').Value = '  +0.85%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.390.73'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.39%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.675'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.66%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.47'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -5.39%  '
$ws.Range('E38').Value = '  -0.32%  '
$ws.Range('E39').Value = '  +0.45%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.965'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +3.05%  '
$ws.Range('E41').Value = '  +10.42%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '82.68'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.58%  '
$ws.Range('E43').Value = '  +1.16%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.42'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.45'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.83%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '6.02'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.70%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0504'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.99%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.965.91'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.04%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '104.55'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.49%  '
$ws.Range('E50').Value = '  -0.10%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0₆0128'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.10%  '
